# Auto-generated edit script: updates cached market-price values in each
# job sheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) to match the latest scrape,
# as produced by the scheduled runner that refreshes this workbook.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2756.1428
$ws.Range("I4").Value = 2658.6
$ws.Range("K4").Value = 2658.6
$ws.Range("M4").Value = -2544.6
$ws.Range("H15").Value = 789.3049
$ws.Range("I15").Value = 789.3049
$ws.Range("K15").Value = 2367.9147
$ws.Range("M15").Value = -2198.9147
$ws.Range("H107").Value = 906.2778
$ws.Range("I107").Value = 761.5833
$ws.Range("K107").Value = 761.5833
$ws.Range("M107").Value = 1158.4167
$ws.Range("H113").Value = 23928
$ws.Range("I113").Value = 41371
$ws.Range("J113").Value = 2124.25
$ws.Range("K113").Value = 41371
$ws.Range("L113").Value = 2124.25
$ws.Range("M113").Value = -38117
$ws.Range("N113").Value = -8632.25
$ws.Range("H131").Value = 3214.1667
$ws.Range("J131").Value = 5333.3335
$ws.Range("L131").Value = 16000.0005
$ws.Range("N131").Value = -26080.0005
$ws.Range("H138").Value = 2378.8
$ws.Range("J138").Value = 2193.383
$ws.Range("L138").Value = 6580.148999999999
$ws.Range("N138").Value = -16860.149
$ws.Range("H141").Value = 2925.8
$ws.Range("I141").Value = 859.6
$ws.Range("J141").Value = 4992
$ws.Range("K141").Value = 2578.8
$ws.Range("L141").Value = 14976
$ws.Range("M141").Value = 2601.2
$ws.Range("N141").Value = -25336

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3823.5122
$ws.Range("I32").Value = 2304.6057
$ws.Range("J32").Value = 13627.363
$ws.Range("K32").Value = 2304.6057
$ws.Range("L32").Value = 13627.363
$ws.Range("M32").Value = -2017.6057
$ws.Range("N32").Value = -14201.363
$ws.Range("H45").Value = 1356.1428
$ws.Range("J45").Value = 1743.5555
$ws.Range("L45").Value = 1743.5555
$ws.Range("N45").Value = -2497.5555
$ws.Range("H61").Value = 73231.75
$ws.Range("I61").Value = 89087
$ws.Range("J61").Value = 25666
$ws.Range("K61").Value = 89087
$ws.Range("L61").Value = 25666
$ws.Range("M61").Value = -88875
$ws.Range("N61").Value = -26090
$ws.Range("H74").Value = 716.23914
$ws.Range("I74").Value = 515.9773
$ws.Range("J74").Value = 5122
$ws.Range("K74").Value = 515.9773
$ws.Range("L74").Value = 5122
$ws.Range("M74").Value = 358.0227
$ws.Range("N74").Value = -6870
$ws.Range("H77").Value = 716.23914
$ws.Range("I77").Value = 515.9773
$ws.Range("J77").Value = 5122
$ws.Range("K77").Value = 2579.8865
$ws.Range("L77").Value = 25610
$ws.Range("M77").Value = 1788.1135
$ws.Range("N77").Value = -34346
$ws.Range("H97").Value = 1793.4286
$ws.Range("I97").Value = 1822.8667
$ws.Range("K97").Value = 1822.8667
$ws.Range("M97").Value = -1326.8667
$ws.Range("H102").Value = 1221.3334
$ws.Range("I102").Value = 1221.3334
$ws.Range("K102").Value = 1221.3334
$ws.Range("M102").Value = 400.6666
$ws.Range("H132").Value = 2473.2727
$ws.Range("I132").Value = 2050.7368
$ws.Range("K132").Value = 6152.2104
$ws.Range("M132").Value = -3622.2104
$ws.Range("H135").Value = 40780
$ws.Range("J135").Value = 40780
$ws.Range("L135").Value = 40780
$ws.Range("N135").Value = -50920
$ws.Range("H136").Value = 73231.75
$ws.Range("I136").Value = 89087
$ws.Range("J136").Value = 25666
$ws.Range("K136").Value = 267261
$ws.Range("L136").Value = 76998
$ws.Range("M136").Value = -264711
$ws.Range("N136").Value = -82098

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 167524.5
$ws.Range("J86").Value = 400684
$ws.Range("L86").Value = 400684
$ws.Range("N86").Value = -402930
$ws.Range("H89").Value = 167524.5
$ws.Range("J89").Value = 400684
$ws.Range("L89").Value = 2003420
$ws.Range("N89").Value = -2014652
$ws.Range("H105").Value = 2087.111
$ws.Range("I105").Value = 2024.0358
$ws.Range("J105").Value = 2307.875
$ws.Range("K105").Value = 2024.0358
$ws.Range("L105").Value = 2307.875
$ws.Range("M105").Value = -277.0358000000001
$ws.Range("N105").Value = -5801.875
$ws.Range("H122").Value = 39666.332
$ws.Range("J122").Value = 39666.332
$ws.Range("L122").Value = 39666.332
$ws.Range("N122").Value = -49466.332

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 860.625
$ws.Range("I16").Value = 840.7143
$ws.Range("K16").Value = 840.7143
$ws.Range("M16").Value = -553.7143
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 1377.7858
$ws.Range("J22").Value = 1822
$ws.Range("L22").Value = 1822
$ws.Range("N22").Value = -2522
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H58").Value = 2175737
$ws.Range("I58").Value = 2719120.5
$ws.Range("K58").Value = 2719120.5
$ws.Range("M58").Value = -2718917.5
$ws.Range("H62").Value = 2879
$ws.Range("I62").Value = 2852
$ws.Range("J62").Value = 2906
$ws.Range("K62").Value = 2852
$ws.Range("L62").Value = 2906
$ws.Range("M62").Value = -2228
$ws.Range("N62").Value = -4154
$ws.Range("H65").Value = 2879
$ws.Range("I65").Value = 2852
$ws.Range("J65").Value = 2906
$ws.Range("K65").Value = 14260
$ws.Range("L65").Value = 14530
$ws.Range("M65").Value = -11140
$ws.Range("N65").Value = -20770
$ws.Range("H113").Value = 860.625
$ws.Range("I113").Value = 840.7143
$ws.Range("K113").Value = 840.7143
$ws.Range("M113").Value = 1329.2857
$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()
$ws.Range("H136").Value = 2175737
$ws.Range("I136").Value = 2719120.5
$ws.Range("K136").Value = 8157361.5
$ws.Range("M136").Value = -8154811.5

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 126124.625
$ws.Range("I11").Value = 200899.6
$ws.Range("K11").Value = 602698.8
$ws.Range("M11").Value = -602558.8
$ws.Range("H49").Value = 3500
$ws.Range("J49").Value = 3500
$ws.Range("L49").Value = 10500
$ws.Range("N49").Value = -10812

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H22").Value = 62504.5
$ws.Range("I22").Value = 45000
$ws.Range("K22").Value = 45000
$ws.Range("M22").Value = -44471
$ws.Range("H102").Value = 2727.8
$ws.Range("I102").Value = 2919.889
$ws.Range("K102").Value = 2919.889
$ws.Range("M102").Value = -1297.889
$ws.Range("H113").Value = 1548.75
$ws.Range("I113").Value = 1197.4
$ws.Range("J113").Value = 1799.7142
$ws.Range("K113").Value = 1197.4
$ws.Range("L113").Value = 1799.7142
$ws.Range("M113").Value = 972.5999999999999
$ws.Range("N113").Value = -6139.7142
$ws.Range("H132").Value = 1204604.9
$ws.Range("I132").Value = 1540885.9
$ws.Range("K132").Value = 4622657.699999999
$ws.Range("M132").Value = -4620127.699999999
$ws.Range("H136").Value = 7378.9414
$ws.Range("J136").Value = 7378.9414
$ws.Range("L136").Value = 22136.8242
$ws.Range("N136").Value = -27236.8242

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3909.889
$ws.Range("J16").Value = 333.66666
$ws.Range("L16").Value = 333.66666
$ws.Range("N16").Value = -673.66666
$ws.Range("H20").Value = 11666.667
$ws.Range("I20").Value = 5000
$ws.Range("K20").Value = 5000
$ws.Range("M20").Value = -4774
$ws.Range("H100").Value = 2999.5
$ws.Range("I100").Value = 2999
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 2999
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -2458
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 2736.3635
$ws.Range("I122").Value = 2677.7778
$ws.Range("K122").Value = 8033.3334
$ws.Range("M122").Value = -5583.3334
$ws.Range("H136").Value = 2683.65
$ws.Range("I136").Value = 2186.353
$ws.Range("K136").Value = 6559.059
$ws.Range("M136").Value = -4009.059

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 887.5833
$ws.Range("I132").Value = 872
$ws.Range("J132").Value = 887.5833
$ws.Range("K132").Value = 2616
$ws.Range("M132").Value = -86
